# Add a new "Save" column (H) to the s_vals sheet, matching the header
# formatting already used by the other header cells (B1:G1), and fill in
# the per-row Save values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header format (bold, bordered, centered - style index 1)
# from G1 onto H1 so the new header cell reuses the same style instead of
# Excel synthesizing a brand-new one.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

$ws.Range("H1").Value = "Save"

$saveValues = @(0, 1, 0, 1, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
